$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need to be forced
# to remain Text (matching the original inlineStr cell type) instead of being
# auto-converted to a numeric value by the COM Value setter.

$ws.Range('D2').Value = '69.880.93'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '3.512.76'
$ws.Range('E3').Value = '  -3.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '193.50'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.29%  '
$ws.Range('D7').Value = '3.502.85'
$ws.Range('E7').Value = '  -3.11%  '
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -6.05%  '
$ws.Range('E11').Value = '  -4.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.64'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.08%  '
$ws.Range('E13').Value = '  -6.01%  '
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('D15').Value = '4.068.64'
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '649.41'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.56%  '
$ws.Range('D17').Value = '69.797.73'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = '3.510.79'
$ws.Range('E18').Value = '  -3.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.44'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.40'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.28%  '
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('E22').Value = '  -4.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.09'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.34'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.46'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.30'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.88%  '
$ws.Range('E27').Value = '  -4.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.11'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.80%  '
$ws.Range('E29').Value = '  -4.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.79'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.99%  '
$ws.Range('E31').Value = '  -8.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.76'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.03%  '
$ws.Range('E34').Value = '  -4.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '557.54'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +8.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.11'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +57.62%  '
$ws.Range('D38').Value = '3.733.66'
$ws.Range('E38').Value = '  -5.11%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.64'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('D41').Value = '0.0₃0792'
$ws.Range('E41').Value = '  -8.32%  '
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.46'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.47%  '
$ws.Range('E46').Value = '  -3.10%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.86'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -5.73%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.35'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.21%  '
$ws.Range('E49').Value = '  -3.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.22'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.75%  '
